# from SSL_CERTFICATE to NORMAL_CERTIFICATE
# Update the two certificate name labels in column A and adjust the
# active-cell selection to A7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Security certificatess"
$ws.Range("A3").Value = "RIPPS certificatess"

$ws.Range("A7").Select()
